$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new (blank) column before the existing "Late" column (old column N),
# pushing "Late" to column O and "Outstanding" to column Q.
$ws.Columns("N").Insert()

# Give the newly inserted column a width similar to its neighbour (column M).
$ws.Columns("N").ColumnWidth = 10.33

# Make "Repayment Schedule" the active sheet/tab and move the selection to R7,
# matching the final view state of the workbook.
$ws.Activate()
$ws.Range("R7").Select()
